$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Force" attribute row above the current row 8 ("Desc"),
# pushing all subsequent rows (LPID_* header info, column headers and
# every PlayerAttNN data row) down by one. Copy row 7's formatting
# (bold font / orange fill / wrap alignment used by the attribute-flag
# rows) so the new row matches its neighbours, then overwrite column A
# with the new attribute name. Columns B:AG stay FALSE, same as the
# copied row.
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Force"
